$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3365.2424
$ws.Range("J138").Value = 3236.818
$ws.Range("L138").Value = 9710.454000000002
$ws.Range("N138").Value = -19990.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1297.75
$ws.Range("I2").Value = 1791
$ws.Range("J2").Value = 1133.3334
$ws.Range("K2").Value = 1791
$ws.Range("L2").Value = 1133.3334
$ws.Range("M2").Value = -1678
$ws.Range("N2").Value = -1359.3334
$ws.Range("H61").Value = 7754322.5
$ws.Range("I61").Value = 18519740
$ws.Range("J61").Value = 3222.04
$ws.Range("K61").Value = 18519740
$ws.Range("L61").Value = 3222.04
$ws.Range("M61").Value = -18519528
$ws.Range("N61").Value = -3646.04
$ws.Range("H74").Value = 932.9697
$ws.Range("I74").Value = 588.7857
$ws.Range("J74").Value = 1186.579
$ws.Range("K74").Value = 588.7857
$ws.Range("L74").Value = 1186.579
$ws.Range("M74").Value = 285.2143
$ws.Range("N74").Value = -2934.579
$ws.Range("H77").Value = 932.9697
$ws.Range("I77").Value = 588.7857
$ws.Range("J77").Value = 1186.579
$ws.Range("K77").Value = 2943.9285
$ws.Range("L77").Value = 5932.895
$ws.Range("M77").Value = 1424.0715
$ws.Range("N77").Value = -14668.895
$ws.Range("H116").Value = 1297.75
$ws.Range("I116").Value = 1791
$ws.Range("J116").Value = 1133.3334
$ws.Range("K116").Value = 1791
$ws.Range("L116").Value = 1133.3334
$ws.Range("M116").Value = 503
$ws.Range("N116").Value = -5721.3334
$ws.Range("H122").Value = 1045.2609
$ws.Range("I122").Value = 955.3158
$ws.Range("J122").Value = 1472.5
$ws.Range("K122").Value = 2865.9474
$ws.Range("L122").Value = 4417.5
$ws.Range("M122").Value = -415.9474
$ws.Range("N122").Value = -9317.5
$ws.Range("H123").Value = 27501.154
$ws.Range("J123").Value = 27501.154
$ws.Range("L123").Value = 27501.154
$ws.Range("N123").Value = -37301.15399999999
$ws.Range("H132").Value = 4865.844
$ws.Range("I132").Value = 4409.8647
$ws.Range("J132").Value = 6974.75
$ws.Range("K132").Value = 13229.5941
$ws.Range("L132").Value = 20924.25
$ws.Range("M132").Value = -10699.5941
$ws.Range("N132").Value = -25984.25
$ws.Range("H136").Value = 7754322.5
$ws.Range("I136").Value = 18519740
$ws.Range("J136").Value = 3222.04
$ws.Range("K136").Value = 55559220
$ws.Range("L136").Value = 9666.119999999999
$ws.Range("M136").Value = -55556670
$ws.Range("N136").Value = -14766.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1297.75
$ws.Range("I3").Value = 1791
$ws.Range("J3").Value = 1133.3334
$ws.Range("K3").Value = 1791
$ws.Range("L3").Value = 1133.3334
$ws.Range("M3").Value = -1677
$ws.Range("N3").Value = -1361.3334
$ws.Range("H80").Value = 1155.1666
$ws.Range("J80").Value = 128.44444
$ws.Range("L80").Value = 128.44444
$ws.Range("N80").Value = -2124.44444
$ws.Range("H83").Value = 1155.1666
$ws.Range("J83").Value = 128.44444
$ws.Range("L83").Value = 642.2221999999999
$ws.Range("N83").Value = -10626.2222
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H134").Value = 3087.2415
$ws.Range("I134").Value = 3555.5386
$ws.Range("J134").Value = 2706.75
$ws.Range("K134").Value = 10666.6158
$ws.Range("L134").Value = 8120.25
$ws.Range("M134").Value = -8131.6158
$ws.Range("N134").Value = -13190.25
$ws.Range("H138").Value = 99780
$ws.Range("J138").Value = 99780
$ws.Range("L138").Value = 99780
$ws.Range("N138").Value = -110060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3729.9348
$ws.Range("I31").Value = 1245.7222
$ws.Range("J31").Value = 5326.9287
$ws.Range("K31").Value = 1245.7222
$ws.Range("L31").Value = 5326.9287
$ws.Range("M31").Value = -950.7221999999999
$ws.Range("N31").Value = -5916.9287
$ws.Range("H34").Value = 3729.9348
$ws.Range("I34").Value = 1245.7222
$ws.Range("J34").Value = 5326.9287
$ws.Range("K34").Value = 1245.7222
$ws.Range("L34").Value = 5326.9287
$ws.Range("M34").Value = -1043.7222
$ws.Range("N34").Value = -5730.9287
$ws.Range("H58").Value = 1986.7407
$ws.Range("I58").Value = 2200.7144
$ws.Range("K58").Value = 2200.7144
$ws.Range("M58").Value = -1997.7144
$ws.Range("H136").Value = 1986.7407
$ws.Range("I136").Value = 2200.7144
$ws.Range("K136").Value = 6602.1432
$ws.Range("M136").Value = -4052.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 90.83871000000001
$ws.Range("J12").Value = 129.75
$ws.Range("L12").Value = 389.25
$ws.Range("N12").Value = -735.25
$ws.Range("H68").Value = 1074.3846
$ws.Range("I68").Value = 1049
$ws.Range("J68").Value = 1078.1177
$ws.Range("K68").Value = 3147
$ws.Range("L68").Value = 3234.3531
$ws.Range("M68").Value = -2336
$ws.Range("N68").Value = -4856.3531
$ws.Range("H71").Value = 1074.3846
$ws.Range("I71").Value = 1049
$ws.Range("J71").Value = 1078.1177
$ws.Range("K71").Value = 9441
$ws.Range("L71").Value = 9703.059300000001
$ws.Range("M71").Value = -5385
$ws.Range("N71").Value = -17815.0593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2187.7778
$ws.Range("I102").Value = 1985.7142
$ws.Range("J102").Value = 2895
$ws.Range("K102").Value = 1985.7142
$ws.Range("L102").Value = 2895
$ws.Range("M102").Value = -363.7141999999999
$ws.Range("N102").Value = -6139
$ws.Range("H122").Value = 1551.65
$ws.Range("I122").Value = 1201.1333
$ws.Range("J122").Value = 2603.2
$ws.Range("K122").Value = 3603.3999
$ws.Range("L122").Value = 7809.599999999999
$ws.Range("M122").Value = -1153.3999
$ws.Range("N122").Value = -12709.6
$ws.Range("H126").Value = 2878.3
$ws.Range("I126").Value = 2828.25
$ws.Range("J126").Value = 2911.6667
$ws.Range("K126").Value = 8484.75
$ws.Range("L126").Value = 8735.000100000001
$ws.Range("M126").Value = -6014.75
$ws.Range("N126").Value = -13675.0001
$ws.Range("H132").Value = 1585.8334
$ws.Range("I132").Value = 780.1177
$ws.Range("J132").Value = 3542.5715
$ws.Range("K132").Value = 2340.3531
$ws.Range("L132").Value = 10627.7145
$ws.Range("M132").Value = 189.6468999999997
$ws.Range("N132").Value = -15687.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 42204.6
$ws.Range("I25").Value = 34503.5
$ws.Range("J25").Value = 47338.668
$ws.Range("K25").Value = 34503.5
$ws.Range("L25").Value = 47338.668
$ws.Range("M25").Value = -34273.5
$ws.Range("N25").Value = -47798.668
$ws.Range("H132").Value = 3073.3918
$ws.Range("I132").Value = 3096.162
$ws.Range("J132").Value = 3050.6216
$ws.Range("K132").Value = 9288.485999999999
$ws.Range("L132").Value = 9151.864799999999
$ws.Range("M132").Value = -6758.485999999999
$ws.Range("N132").Value = -14211.8648

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 30550
$ws.Range("J92").Value = 30550
$ws.Range("L92").Value = 30550
$ws.Range("N92").Value = -35542
$ws.Range("H100").Value = 515.75
$ws.Range("I100").Value = 302.53333
$ws.Range("J100").Value = 1155.4
$ws.Range("K100").Value = 605.06666
$ws.Range("L100").Value = 2310.8
$ws.Range("M100").Value = -64.06665999999996
$ws.Range("N100").Value = -3392.8
$ws.Range("H122").Value = 1718.4375
$ws.Range("I122").Value = 1499.091
$ws.Range("K122").Value = 4497.272999999999
$ws.Range("M122").Value = -2047.272999999999
$ws.Range("H123").Value = 24610.059
$ws.Range("J123").Value = 24610.059
$ws.Range("L123").Value = 24610.059
$ws.Range("N123").Value = -34410.059
$ws.Range("H124").Value = 27981.857
$ws.Range("J124").Value = 27981.857
$ws.Range("L124").Value = 27981.857
$ws.Range("N124").Value = -37801.857
